$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder rows 2-4 (login story now sorted after the two priority-0 stories) ---
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "Nurse"
$ws.Cells.Item(2, 3).Value = "I want to create a new patient record with their personal details"
$ws.Cells.Item(2, 4).Value = "so that I can add new patients to the system"

$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = "Nurse"
$ws.Cells.Item(3, 3).Value = "I want to schedule an appointment without allowing double-booking for the same doctor"
$ws.Cells.Item(3, 4).Value = "so that a patient can see a doctor"

$ws.Cells.Item(4, 1).Value = 0.5
$ws.Cells.Item(4, 2).Value = "Nurse & Administrator"
$ws.Cells.Item(4, 3).Value = "I wish to log in using my credentials."
$ws.Cells.Item(4, 4).Value = "so that I can securely access the system"

# --- Insert 3 new rows for the missing user stories (pushes old rows 9-11 to 12-14) ---
$ws.Range("A9:D11").EntireRow.Insert()

$ws.Cells.Item(9, 1).Value = 1
$ws.Cells.Item(9, 2).Value = "Nurse"
$ws.Cells.Item(9, 3).Value = "I want to edit patient information"
$ws.Cells.Item(9, 4).Value = "so that I can update patient information"

$ws.Cells.Item(10, 1).Value = 1
$ws.Cells.Item(10, 2).Value = "Nurse"
$ws.Cells.Item(10, 3).Value = "I want to search for patient visit information"
$ws.Cells.Item(10, 4).Value = "so that I can see patient visit information"

$ws.Cells.Item(11, 1).Value = 1
$ws.Cells.Item(11, 2).Value = "Nurse"
$ws.Cells.Item(11, 3).Value = "I want to search for patient appointment information"
$ws.Cells.Item(11, 4).Value = "so that I can see patient appointment information"

# --- Update sort state / selection to reflect the new data range ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A14"))
$ws.Sort.SetRange($ws.Range("A2:D14"))
$ws.Sort.Header = $false
[void]$ws.Sort.Apply()

[void]$ws.Range("C18").Select()
